$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Credentials")

# New credentials text for Kode Cloud (C4), matching the style used in C3
$credText = "id : Sandhyaramanipradeep@yahoo.com" + [char]10 + "pass : Happy123!"

# Set the cell value (creates a new shared string entry)
$ws.Range("C4").Value = $credText

# Match style of C3 (wrap text) and row height (2 lines tall)
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 29

# Widen column C to fit the new, longer content
$ws.Columns.Item(3).ColumnWidth = 43.3
